# Weekly update: insert two new price rows for
# "Feria Lagunitas de Puerto Montt - Repollo" (rows 818-819), pushing the
# existing rows 818-875 down to 820-877.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 818, shifting everything currently at 818:875
# down to 820:877 (and extending the used range to row 877).
$ws.Range("A818:R819").Insert()

# New row 818
$ws.Range("A818").Value = 4
$ws.Range("B818").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C818").Value = "Los Lagos"
$ws.Range("D818").Value = 45265
$ws.Range("E818").Value = 10
$ws.Range("F818").Value = 100112006
$ws.Range("G818").Value = "Repollo"
$ws.Range("H818").Value = "Copenhague"
$ws.Range("I818").Value = "Primera"
$ws.Range("J818").Value = 600
$ws.Range("K818").Value = 1800
$ws.Range("L818").Value = 2000
$ws.Range("M818").Value = 1900
$ws.Range("N818").Value = "$/unidad"
$ws.Range("O818").Value = "Región Metropolitana"
$ws.Range("P818").Value = 1900
$ws.Range("Q818").Value = 1
$ws.Range("R818").Value = "Hortaliza"

# New row 819
$ws.Range("A819").Value = 4
$ws.Range("B819").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C819").Value = "Los Lagos"
$ws.Range("D819").Value = 45265
$ws.Range("E819").Value = 10
$ws.Range("F819").Value = 100112006
$ws.Range("G819").Value = "Repollo"
$ws.Range("H819").Value = "Crespo record"
$ws.Range("I819").Value = "Primera"
$ws.Range("J819").Value = 800
$ws.Range("K819").Value = 1500
$ws.Range("L819").Value = 1500
$ws.Range("M819").Value = 1500
$ws.Range("N819").Value = "$/unidad"
$ws.Range("O819").Value = "Región Metropolitana"
$ws.Range("P819").Value = 1500
$ws.Range("Q819").Value = 1
$ws.Range("R819").Value = "Hortaliza"
